$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 54-55, shifting existing rows 54-135 down to 56-137
$ws.Rows("54:55").Insert()

# Populate new row 54 (Lapins / Primera)
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = 'Macroferia Regional de Talca'
$ws.Range("C54").Value = 'Maule'
$ws.Range("D54").Value = 44546
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = 'Frutos de hueso (carozo)'
$ws.Range("I54").Value = 100103001
$ws.Range("J54").Value = 'Cereza'
$ws.Range("K54").Value = 'Lapins'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 200
$ws.Range("N54").Value = 6000
$ws.Range("O54").Value = 6000
$ws.Range("P54").Value = 6000
$ws.Range("Q54").Value = '$/bandeja 10 kilos'
$ws.Range("R54").Value = 'Provincia de Curicó'
$ws.Range("S54").Value = 600
$ws.Range("T54").Value = 10

# Populate new row 55 (Rainier / Primera)
$ws.Range("A55").Value = 5
$ws.Range("B55").Value = 'Macroferia Regional de Talca'
$ws.Range("C55").Value = 'Maule'
$ws.Range("D55").Value = 44546
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 'Fruta'
$ws.Range("G55").Value = 100103
$ws.Range("H55").Value = 'Frutos de hueso (carozo)'
$ws.Range("I55").Value = 100103001
$ws.Range("J55").Value = 'Cereza'
$ws.Range("K55").Value = 'Rainier'
$ws.Range("L55").Value = 'Primera'
$ws.Range("M55").Value = 180
$ws.Range("N55").Value = 8000
$ws.Range("O55").Value = 8000
$ws.Range("P55").Value = 8000
$ws.Range("Q55").Value = '$/bandeja 10 kilos'
$ws.Range("R55").Value = 'Provincia de Curicó'
$ws.Range("S55").Value = 800
$ws.Range("T55").Value = 10

